# Removed dependency of dose multiplier: delete the "Dose multiplier
# interval" column (column G) entirely, shifting On/off (old H) and
# Dose vol. (old I) one column to the left.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("G").Select() | Out-Null
$ws.Columns("G").Delete() | Out-Null
